{"js": "// Update the worksheet date and the 25 division-problem \"answers\" cells\n// to the next day's generated set, per the commit:\n//   \"Update master to output generated at 4250d90\"\n//\n// Each (old, new) pair below is an exact, unique full-text match of a\n// single w:t run in the document (the title date paragraph, plus the\n// 25 populated table cells out of the 6x5 answer grid). We locate each\n// old string with body.search(...) and overwrite the whole run's text\n// via Range.insertText(new, \"Replace\").\n\nconst replacements = [\n  [\"2024-08-06 Tuesday\", \"2024-08-07 Wednesday\"],\n  [\"32\u00f79=3, 5\", \"36\u00f74=9, 0\"],\n  [\"62\u00f77=8, 6\", \"38\u00f72=19, 0\"],\n  [\"50\u00f79=5, 5\", \"12\u00f72=6, 0\"],\n  [\"12\u00f77=1, 5\", \"91\u00f77=13, 0\"],\n  [\"17\u00f77=2, 3\", \"90\u00f79=10, 0\"],\n  [\"44\u00f76=7, 2\", \"64\u00f74=16, 0\"],\n  [\"19\u00f77=2, 5\", \"62\u00f74=15, 2\"],\n  [\"21\u00f79=2, 3\", \"80\u00f73=26, 2\"],\n  [\"97\u00f73=32, 1\", \"77\u00f78=9, 5\"],\n  [\"30\u00f74=7, 2\", \"34\u00f77=4, 6\"],\n  [\"77\u00f72=38, 1\", \"25\u00f75=5, 0\"],\n  [\"16\u00f76=2, 4\", \"58\u00f73=19, 1\"],\n  [\"88\u00f75=17, 3\", \"68\u00f73=22, 2\"],\n  [\"84\u00f74=21, 0\", \"79\u00f72=39, 1\"],\n  [\"26\u00f72=13, 0\", \"89\u00f72=44, 1\"],\n  [\"24\u00f72=12, 0\", \"10\u00f78=1, 2\"],\n  [\"45\u00f77=6, 3\", \"24\u00f74=6, 0\"],\n  [\"46\u00f77=6, 4\", \"24\u00f79=2, 6\"],\n  [\"64\u00f76=10, 4\", \"28\u00f79=3, 1\"],\n  [\"79\u00f76=13, 1\", \"48\u00f75=9, 3\"],\n  [\"96\u00f72=48, 0\", \"14\u00f75=2, 4\"],\n  [\"50\u00f78=6, 2\", \"25\u00f75=5, 0\"],\n  [\"71\u00f74=17, 3\", \"16\u00f78=2, 0\"],\n  [\"54\u00f77=7, 5\", \"28\u00f74=7, 0\"],\n  [\"52\u00f73=17, 1\", \"73\u00f78=9, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division-problem \"answers\" cells\n# to the next day's generated set, per the commit:\n#   \"Update master to output generated at 4250d90\"\n#\n# Each (old, new) pair is the exact, unique full text of a single run\n# (the title date paragraph, plus the 25 populated cells of the 6x5\n# answer grid). We drive Word's Find/Replace across the whole story so\n# every occurrence (exactly one, per uniqueness of each string) is\n# swapped in place, preserving all formatting/structure.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-06 Tuesday\", \"2024-08-07 Wednesday\"),\n    @(\"32\u00f79=3, 5\", \"36\u00f74=9, 0\"),\n    @(\"62\u00f77=8, 6\", \"38\u00f72=19, 0\"),\n    @(\"50\u00f79=5, 5\", \"12\u00f72=6, 0\"),\n    @(\"12\u00f77=1, 5\", \"91\u00f77=13, 0\"),\n    @(\"17\u00f77=2, 3\", \"90\u00f79=10, 0\"),\n    @(\"44\u00f76=7, 2\", \"64\u00f74=16, 0\"),\n    @(\"19\u00f77=2, 5\", \"62\u00f74=15, 2\"),\n    @(\"21\u00f79=2, 3\", \"80\u00f73=26, 2\"),\n    @(\"97\u00f73=32, 1\", \"77\u00f78=9, 5\"),\n    @(\"30\u00f74=7, 2\", \"34\u00f77=4, 6\"),\n    @(\"77\u00f72=38, 1\", \"25\u00f75=5, 0\"),\n    @(\"16\u00f76=2, 4\", \"58\u00f73=19, 1\"),\n    @(\"88\u00f75=17, 3\", \"68\u00f73=22, 2\"),\n    @(\"84\u00f74=21, 0\", \"79\u00f72=39, 1\"),\n    @(\"26\u00f72=13, 0\", \"89\u00f72=44, 1\"),\n    @(\"24\u00f72=12, 0\", \"10\u00f78=1, 2\"),\n    @(\"45\u00f77=6, 3\", \"24\u00f74=6, 0\"),\n    @(\"46\u00f77=6, 4\", \"24\u00f79=2, 6\"),\n    @(\"64\u00f76=10, 4\", \"28\u00f79=3, 1\"),\n    @(\"79\u00f76=13, 1\", \"48\u00f75=9, 3\"),\n    @(\"96\u00f72=48, 0\", \"14\u00f75=2, 4\"),\n    @(\"50\u00f78=6, 2\", \"25\u00f75=5, 0\"),\n    @(\"71\u00f74=17, 3\", \"16\u00f78=2, 0\"),\n    @(\"54\u00f77=7, 5\", \"28\u00f74=7, 0\"),\n    @(\"52\u00f73=17, 1\", \"73\u00f78=9, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)  # 2 = wdReplaceAll\n\n    if (-not $ok) {\n        throw \"Could not find text to replace: '$oldText'\"\n    }\n}\n"}
